$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Devin Vassell, SG,SF, San Antonio Spurs -> Trae Young, PG, Atlanta Hawks
$ws.Range("A5").Value = "Trae Young"
$ws.Range("B5").Value = "PG"
$ws.Range("C5").Value = "Atlanta Hawks"

# Row 8: Kelly Oubre Jr., SG,SF, Philadelphia 76ers -> Devin Vassell, SG,SF, San Antonio Spurs
$ws.Range("A8").Value = "Devin Vassell"
$ws.Range("C8").Value = "San Antonio Spurs"

# Row 10: Alperen Sengun, C, Houston Rockets -> Walker Kessler, C, Utah Jazz
$ws.Range("A10").Value = "Walker Kessler"
$ws.Range("C10").Value = "Utah Jazz"

# Row 11: Walker Kessler, C, Utah Jazz -> Myles Turner, C, Indiana Pacers
$ws.Range("A11").Value = "Myles Turner"
$ws.Range("C11").Value = "Indiana Pacers"

# Row 12: Trae Young, PG, Atlanta Hawks -> LeBron James, SF,PF, Los Angeles Lakers
$ws.Range("A12").Value = "LeBron James"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "Los Angeles Lakers"

# Row 13: LeBron James, SF,PF, Los Angeles Lakers -> Jamal Murray, PG,SG, Denver Nuggets
$ws.Range("A13").Value = "Jamal Murray"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Denver Nuggets"

# Row 14: Amen Thompson, SG,SF, Houston Rockets -> Kelly Oubre Jr., SG,SF, Philadelphia 76ers
$ws.Range("A14").Value = "Kelly Oubre Jr."
$ws.Range("C14").Value = "Philadelphia 76ers"

# Row 16: Caris LeVert, SG,SF, Cleveland Cavaliers -> Klay Thompson, SG,SF, Dallas Mavericks
$ws.Range("A16").Value = "Klay Thompson"
$ws.Range("C16").Value = "Dallas Mavericks"
